# Auto-generated: apply scheduled-runner market data refresh
# Updates columns H-N (computed market price / profit fields) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 361
$ws.Range("I2").Value = 146
$ws.Range("J2").Value = 898.5
$ws.Range("K2").Value = 146
$ws.Range("L2").Value = 898.5
$ws.Range("M2").Value = -33
$ws.Range("N2").Value = -1124.5
# row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 677.12
$ws.Range("I28").Value = 597.8946999999999
$ws.Range("K28").Value = 597.8946999999999
$ws.Range("M28").Value = -112.8946999999999
# row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 10845.875
$ws.Range("I137").Value = 1523.3182
$ws.Range("J137").Value = 31355.5
$ws.Range("K137").Value = 4569.9546
$ws.Range("L137").Value = 94066.5
$ws.Range("M137").Value = -2019.9546
$ws.Range("N137").Value = -99166.5
# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 1626060.9
$ws.Range("I138").Value = 1909.8
$ws.Range("J138").Value = 2364311.2
$ws.Range("K138").Value = 5729.4
$ws.Range("L138").Value = 7092933.600000001
$ws.Range("M138").Value = -589.3999999999996
$ws.Range("N138").Value = -7103213.600000001

$ws = $wb.Worksheets.Item("ARM")
# row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 930.2
$ws.Range("I2").Value = 912.875
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 912.875
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -799.875
$ws.Range("N2").Value = -1225.5
# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 25887.934
$ws.Range("I32").Value = 29398.447
$ws.Range("J32").Value = 6830.857
$ws.Range("K32").Value = 29398.447
$ws.Range("L32").Value = 6830.857
$ws.Range("M32").Value = -29111.447
$ws.Range("N32").Value = -7404.857
# row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 3728.7
$ws.Range("I45").Value = 1493
$ws.Range("K45").Value = 1493
$ws.Range("M45").Value = -1116
# row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 5913.4316
$ws.Range("I61").Value = 3708.647
$ws.Range("K61").Value = 3708.647
$ws.Range("M61").Value = -3496.647
# row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 21767.426
$ws.Range("I110").Value = 24388.4
$ws.Range("K110").Value = 24388.4
$ws.Range("M110").Value = -22343.4
# row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 930.2
$ws.Range("I116").Value = 912.875
$ws.Range("J116").Value = 999.5
$ws.Range("K116").Value = 912.875
$ws.Range("L116").Value = 999.5
$ws.Range("M116").Value = 1381.125
$ws.Range("N116").Value = -5587.5
# row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2488.4707
$ws.Range("I122").Value = 1805.091
$ws.Range("K122").Value = 5415.272999999999
$ws.Range("M122").Value = -2965.272999999999
# row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3857.4102
$ws.Range("I132").Value = 3915.0344
$ws.Range("J132").Value = 3690.3
$ws.Range("K132").Value = 11745.1032
$ws.Range("L132").Value = 11070.9
$ws.Range("M132").Value = -9215.1032
$ws.Range("N132").Value = -16130.9
# row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 5913.4316
$ws.Range("I136").Value = 3708.647
$ws.Range("K136").Value = 11125.941
$ws.Range("M136").Value = -8575.940999999999

$ws = $wb.Worksheets.Item("BSM")
# row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 930.2
$ws.Range("I3").Value = 912.875
$ws.Range("J3").Value = 999.5
$ws.Range("K3").Value = 912.875
$ws.Range("L3").Value = 999.5
$ws.Range("M3").Value = -798.875
$ws.Range("N3").Value = -1227.5
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 1925.45
$ws.Range("I86").Value = 2194.375
$ws.Range("J86").Value = 849.75
$ws.Range("K86").Value = 2194.375
$ws.Range("L86").Value = 849.75
$ws.Range("M86").Value = -1071.375
$ws.Range("N86").Value = -3095.75
# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1925.45
$ws.Range("I89").Value = 2194.375
$ws.Range("J89").Value = 849.75
$ws.Range("K89").Value = 10971.875
$ws.Range("L89").Value = 4248.75
$ws.Range("M89").Value = -5355.875
$ws.Range("N89").Value = -15480.75
# row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 4911
$ws.Range("I105").Value = 6141
$ws.Range("J105").Value = 4501
$ws.Range("K105").Value = 6141
$ws.Range("L105").Value = 4501
$ws.Range("M105").Value = -4394
$ws.Range("N105").Value = -7995
# row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 7501.385
$ws.Range("I134").Value = 8891.111000000001
$ws.Range("J134").Value = 4374.5
$ws.Range("K134").Value = 26673.333
$ws.Range("L134").Value = 13123.5
$ws.Range("M134").Value = -24138.333
$ws.Range("N134").Value = -18193.5

$ws = $wb.Worksheets.Item("CRP")
# row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 3518.9092
$ws.Range("I105").Value = 3285.2856
$ws.Range("K105").Value = 3285.2856
$ws.Range("M105").Value = -1538.2856

$ws = $wb.Worksheets.Item("CUL")
# row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 752.1667
$ws.Range("I5").Value = 721.2222
$ws.Range("K5").Value = 2163.6666
$ws.Range("M5").Value = -2051.6666
# row 63 (Leve Item ID 12866)
$ws.Range("H63").Value = 12414
$ws.Range("J63").Value = 4896.8
$ws.Range("L63").Value = 14690.4
$ws.Range("N63").Value = -16188.4
# row 66 (Leve Item ID 12866)
$ws.Range("H66").Value = 12414
$ws.Range("J66").Value = 4896.8
$ws.Range("L66").Value = 44071.2
$ws.Range("N66").Value = -51559.2
# row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 752.1667
$ws.Range("I135").Value = 721.2222
$ws.Range("K135").Value = 6490.999800000001
$ws.Range("M135").Value = -3955.999800000001
# row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 2126.818
$ws.Range("J137").Value = 2882.6667
$ws.Range("L137").Value = 8648.000100000001
$ws.Range("N137").Value = -18848.0001

$ws = $wb.Worksheets.Item("GSM")
# row 49 (Leve Item ID 4232)
$ws.Range("H49").Value = 30024
$ws.Range("J49").Value = 30024
$ws.Range("L49").Value = 30024
$ws.Range("N49").Value = -30392
# row 62 (Leve Item ID 11983)
$ws.Range("H62").Value = 59996.715
$ws.Range("I62").Value = 59988
$ws.Range("K62").Value = 59988
$ws.Range("M62").Value = -59302
# row 65 (Leve Item ID 11983)
$ws.Range("H65").Value = 59996.715
$ws.Range("I65").Value = 59988
$ws.Range("K65").Value = 179964
$ws.Range("M65").Value = -176532
# row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 13999.5
$ws.Range("I80").Value = 18999.5
$ws.Range("K80").Value = 18999.5
$ws.Range("M80").Value = -18001.5
# row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 13999.5
$ws.Range("I83").Value = 18999.5
$ws.Range("K83").Value = 94997.5
$ws.Range("M83").Value = -90005.5
# row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1683.3684
$ws.Range("I102").Value = 1469.7059
$ws.Range("K102").Value = 1469.7059
$ws.Range("M102").Value = 152.2941000000001
# row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 2430.3635
$ws.Range("I113").Value = 2192.7778
$ws.Range("J113").Value = 3499.5
$ws.Range("K113").Value = 2192.7778
$ws.Range("L113").Value = 3499.5
$ws.Range("M113").Value = -22.77779999999984
$ws.Range("N113").Value = -7839.5
# row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 5574.2573
$ws.Range("I132").Value = 6808.2856
$ws.Range("K132").Value = 20424.8568
$ws.Range("M132").Value = -17894.8568

$ws = $wb.Worksheets.Item("LTW")
# row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 986.3570999999999
$ws.Range("J16").Value = 1100.25
$ws.Range("L16").Value = 1100.25
$ws.Range("N16").Value = -1440.25
# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2025
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 2200
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 2200
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -2790
# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2025
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 2200
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 2200
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -2414
# row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 3160.28
$ws.Range("I40").Value = 3043.7827
$ws.Range("K40").Value = 3043.7827
$ws.Range("M40").Value = -2907.7827
# row 42 (Leve Item ID 4333)
$ws.Range("H42").Value = 16363.637
$ws.Range("J42").Value = 16363.637
$ws.Range("L42").Value = 16363.637
$ws.Range("N42").Value = -17489.637
# row 49 (Leve Item ID 4333)
$ws.Range("H49").Value = 16363.637
$ws.Range("J49").Value = 16363.637
$ws.Range("L49").Value = 16363.637
$ws.Range("N49").Value = -16657.637
# row 74 (Leve Item ID 11990)
$ws.Range("H74").Value = 126874.125
$ws.Range("I74").Value = 116248.75
$ws.Range("K74").Value = 116248.75
$ws.Range("M74").Value = -115250.75
# row 77 (Leve Item ID 11990)
$ws.Range("H77").Value = 126874.125
$ws.Range("I77").Value = 116248.75
$ws.Range("K77").Value = 348746.25
$ws.Range("M77").Value = -343754.25
# row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2465.6667
$ws.Range("I122").Value = 2299.8572
$ws.Range("K122").Value = 6899.571599999999
$ws.Range("M122").Value = -4449.571599999999
# row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 4155.3213
$ws.Range("I136").Value = 2814.1875
$ws.Range("K136").Value = 8442.5625
$ws.Range("M136").Value = -5892.5625

$ws = $wb.Worksheets.Item("WVR")
# row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 11926.571
$ws.Range("I62").Value = 11426.714
$ws.Range("K62").Value = 11426.714
$ws.Range("M62").Value = -10802.714
# row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 11926.571
$ws.Range("I65").Value = 11426.714
$ws.Range("K65").Value = 57133.57
$ws.Range("M65").Value = -54013.57
# row 80 (Leve Item ID 10911)
$ws.Range("H80").Value = 90301
$ws.Range("J80").Value = 90301
$ws.Range("L80").Value = 90301
$ws.Range("N80").Value = -92297
# row 83 (Leve Item ID 10911)
$ws.Range("H83").Value = 90301
$ws.Range("J83").Value = 90301
$ws.Range("L83").Value = 270903
$ws.Range("N83").Value = -280887
# row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 15628872
$ws.Range("I122").Value = 21743158
$ws.Range("K122").Value = 65229474
$ws.Range("M122").Value = -65227024
# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 13198026
$ws.Range("I132").Value = 19287536
$ws.Range("J132").Value = 4090.3333
$ws.Range("K132").Value = 57862608
$ws.Range("L132").Value = 12270.9999
$ws.Range("M132").Value = -57860078
$ws.Range("N132").Value = -17330.9999
